# Commit: "added tabview for search results, todo: fix recursion error"
#
# Semantic change (per the OOXML diff):
#   - New row 8 added to Sheet1 with F8 = "go" (shared string already
#     present in sharedStrings.xml at index 12, reused -> sst count 72->73,
#     uniqueCount unchanged at 71).
#   - Worksheet dimension grows from A1:L6 to A1:L8.
#   - Selection/active cell moves from L6 to G8.
#
# (The incidental Office-chrome metadata in the diff -- absPath url,
# xr:revisionPtr GUIDs, workbookView window geometry, and the empty
# xl/persons/person.xml stub -- are host/session artifacts Excel stamps on
# save and aren't reachable/meaningful through the object model, so they're
# intentionally left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new cell value; matches shared string "go" already used by A2.
$ws.Range("F8").Value = "go"

# Move the selection/active cell to G8, matching the saved sheetView state.
$ws.Range("G8").Select() | Out-Null
